$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cryptocurrency price/volume data (and a few name/link cells
# that shifted because a new "Frax" row was inserted at row 44).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.074.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.54%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.99%  '

# Row 4
$ws.Range("E4").Value = '  +0.31%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.38%  '

# Row 6
$ws.Range("E6").Value = '  +0.19%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3841'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.29%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08608'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.09%  '

# Row 10
$ws.Range("E10").Value = '  -2.38%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.49'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.14%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.309'
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.25%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.878.22'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.78%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.192'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.84%  '

# Row 16
$ws.Range("E16").Value = '  +0.34%  '

# Row 17
$ws.Range("E17").Value = '  -2.30%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.75%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06625'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.082'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.42%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.112.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.54%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.264'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.59%  '

# Row 26
$ws.Range("E26").Value = '  -0.95%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.099.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.25%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.96%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.79%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1054'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.060'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.601'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.98%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.597'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.44%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.626'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.09%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02441'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.53%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06580'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.57%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2175'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.63%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.211'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.98%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6375'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.36%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.66%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.890'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.64%  '

# Row 44
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.23%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.37%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5998'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.01%  '

# Row 47
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.283'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.32%  '

# Row 48
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.669'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.59%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.991'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.68%  '

# Row 50
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.224'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.89%  '

# Row 51
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.32%  '
